# Run the Jewels RTJ3 tracklist update (replaces old Dave - Psychodrama data)
$wb = $excel.ActiveWorkbook

$trackData = @(
    @("B1", "Title/Composer"),
    @("D1", "Performer"),
    @("E1", "Time"),
    @("A2", 1),
    @("B2", "Down"),
    @("C2", "El-P / Killer Mike"),
    @("D2", "Run the Jewels feat. Joi"),
    @("E2", 0.1451388888888889),
    @("A3", 2),
    @("B3", "Talk to Me"),
    @("C3", "El-P / Killer Mike"),
    @("D3", "Run the Jewels"),
    @("E3", 0.10486111111111111),
    @("A4", 3),
    @("B4", "Legend Has It"),
    @("C4", "El-P / Killer Mike"),
    @("D4", "Run the Jewels"),
    @("E4", 0.1423611111111111),
    @("A5", 4),
    @("B5", "Call Ticketron"),
    @("C5", "El-P / Killer Mike"),
    @("D5", "Run the Jewels"),
    @("E5", 0.13749999999999998),
    @("A6", 5),
    @("B6", "Hey Kids"),
    @("C6", "Danny Brown / El-P / Killer Mike"),
    @("D6", "Run the Jewels feat. Danny Brown"),
    @("E6", 0.1326388888888889),
    @("A7", 6),
    @("B7", "Stay Gold"),
    @("C7", "El-P / Killer Mike"),
    @("D7", "Run the Jewels"),
    @("E7", 0.14375000000000002),
    @("A8", 7),
    @("B8", "Don't Get Captured"),
    @("C8", "El-P / Killer Mike"),
    @("D8", "Run the Jewels"),
    @("E8", 0.13333333333333333),
    @("A9", 8),
    @("B9", "Thieves! (Screamed the Ghost)"),
    @("C9", "Boots / El-P / Killer Mike"),
    @("D9", "Run the Jewels feat. Tande Adebimpe"),
    @("E9", 0.16805555555555554),
    @("A10", 9),
    @("B10", 2100),
    @("C10", "Boots / El-P / Killer Mike"),
    @("D10", "Run the Jewels feat. Boots"),
    @("E10", 0.1673611111111111),
    @("A11", 10),
    @("B11", "Panther Like a Panther [Miracle Mix]"),
    @("C11", "Boots / El-P / Killer Mike"),
    @("D11", "Run the Jewels feat. Trina"),
    @("E11", 0.15347222222222223),
    @("A12", 11),
    @("B12", "Everybody Stay Calm"),
    @("C12", "El-P / Killer Mike"),
    @("D12", "Run the Jewels"),
    @("E12", 0.12361111111111112),
    @("A13", 12),
    @("B13", "Oh Mama"),
    @("C13", "El-P / Killer Mike"),
    @("D13", "Run the Jewels"),
    @("E13", 0.15),
    @("A14", 13),
    @("B14", "Thursday in the Danger Room"),
    @("C14", "El-P / Killer Mike"),
    @("D14", "Run the Jewels feat: Kamasi Washington"),
    @("E14", 0.18194444444444444),
    @("A15", 14),
    @("B15", "A Report to Your Shareholders/Kill Your Masters"),
    @("C15", "El-P / Killer Mike / Zack de la Rocha"),
    @("D15", "Run the Jewels"),
    @("E15", 0.25972222222222224)
)

$dataSheets = @($wb.Worksheets.Item("Sheet1"), $wb.Worksheets.Item("Sheet3"))

foreach ($ws in $dataSheets) {
    # Row 1 header no longer carries a leading index number or a separate
    # "Composer" header cell - clear those two legacy cells first.
    $ws.Range("A1").ClearContents()
    $ws.Range("C1").ClearContents()

    foreach ($item in $trackData) {
        $ws.Range($item[0]).Value = $item[1]
    }
}

# Defined names: "dave1" -> "runthejewels3", range extended from row 12 to row 15.
# (Renaming Name objects in place is unreliable in this host when there are
# two same-named local scopes, so drop and recreate them instead.)
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}
$ws1.Names.Add("runthejewels3", "=Sheet1!`$A`$1:`$E`$15")
$ws3.Names.Add("runthejewels3", "=Sheet3!`$A`$1:`$E`$15")

# Sheet1 / Sheet3 selection: mirror the new used range
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("A1:E15").Select()

# Sheet2 is fully formula-driven off Sheet1 and recalculates automatically;
# just move its reported selection down to match the extra rows.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A3:K18").Select()

$ws2.Activate()
